$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "User ID"
$ws.Range("C1").Value = "Product ID"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "Total Price"

# Best-fit column widths (mirrors Excel's "AutoFit selection" behaviour for
# the header row) expressed in character-width units.
$ws.Columns.Item(1).ColumnWidth = 2.1666666666666665
$ws.Columns.Item(2).ColumnWidth = 6.666666666666667
$ws.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 8.0
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
